# Scheduled-runner price/profit refresh: updates the currentAveragePrice*,
# LevePriceNQ/HQ and LeveProfitNQ/HQ columns (H-N) on a handful of rows
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets with freshly pulled
# market-board figures. Leve metadata columns (A-G) are left untouched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 31.625
$ws.Range("I6").Value = 31.625
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 94.875
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 17.125
$ws.Range("N6").ClearContents()
$ws.Range("H8").Value = 1452
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 1452
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 4356
$ws.Range("M8").ClearContents()
$ws.Range("N8").Value = -4634
$ws.Range("H17").Value = 1418.4562
$ws.Range("J17").Value = 1408.0714
$ws.Range("L17").Value = 4224.2142
$ws.Range("N17").Value = -4560.2142
$ws.Range("H31").Value = 7134.4
$ws.Range("I31").Value = 8793.25
$ws.Range("K31").Value = 26379.75
$ws.Range("M31").Value = -26149.75
$ws.Range("H33").Value = 735.65625
$ws.Range("I33").Value = 814.6
$ws.Range("J33").Value = 453.7143
$ws.Range("K33").Value = 814.6
$ws.Range("L33").Value = 453.7143
$ws.Range("M33").Value = -585.6
$ws.Range("N33").Value = -911.7143
$ws.Range("H52").Value = 300
$ws.Range("I52").Value = 300
$ws.Range("K52").Value = 900
$ws.Range("M52").Value = -740
$ws.Range("H54").Value = 15737.5
$ws.Range("I54").Value = 16975
$ws.Range("K54").Value = 16975
$ws.Range("M54").Value = -16489
$ws.Range("H59").Value = 949.6667
$ws.Range("I59").Value = 850
$ws.Range("K59").Value = 2550
$ws.Range("M59").Value = -1993
$ws.Range("H69").Value = 13746.091
$ws.Range("I69").Value = 9702.166999999999
$ws.Range("K69").Value = 29106.501
$ws.Range("M69").Value = -28232.501
$ws.Range("H72").Value = 13746.091
$ws.Range("I72").Value = 9702.166999999999
$ws.Range("K72").Value = 87319.503
$ws.Range("M72").Value = -82951.503
$ws.Range("H74").Value = 5130.385
$ws.Range("I74").Value = 4969.5
$ws.Range("K74").Value = 4969.5
$ws.Range("M74").Value = -4033.5
$ws.Range("H77").Value = 5130.385
$ws.Range("I77").Value = 4969.5
$ws.Range("K77").Value = 24847.5
$ws.Range("M77").Value = -20167.5
$ws.Range("H96").Value = 1569.4828
$ws.Range("I96").Value = 1955.6666
$ws.Range("K96").Value = 5866.9998
$ws.Range("M96").Value = -4493.9998
$ws.Range("H100").Value = 1112.4117
$ws.Range("I100").Value = 1142.1111
$ws.Range("K100").Value = 1142.1111
$ws.Range("M100").Value = -601.1111000000001
$ws.Range("H111").Value = 2301.6428
$ws.Range("I111").Value = 1709.4615
$ws.Range("K111").Value = 5128.3845
$ws.Range("M111").Value = -2061.3845
$ws.Range("H113").Value = 4543.52
$ws.Range("I113").Value = 3924.5
$ws.Range("K113").Value = 3924.5
$ws.Range("M113").Value = -670.5
$ws.Range("H116").Value = 3034.647
$ws.Range("I116").Value = 2639.3333
$ws.Range("K116").Value = 2639.3333
$ws.Range("M116").Value = 802.6667000000002
$ws.Range("H125").Value = 2451.8
$ws.Range("I125").Value = 805
$ws.Range("K125").Value = 7245
$ws.Range("M125").Value = -4785

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5048.845
$ws.Range("I32").Value = 4449.543
$ws.Range("K32").Value = 4449.543
$ws.Range("M32").Value = -4162.543
$ws.Range("H45").Value = 10135.7
$ws.Range("I45").Value = 9924.294
$ws.Range("K45").Value = 9924.294
$ws.Range("M45").Value = -9547.294
$ws.Range("H63").Value = 3565.8667
$ws.Range("I63").Value = 3824.5
$ws.Range("J63").Value = 3547.3928
$ws.Range("K63").Value = 3824.5
$ws.Range("L63").Value = 3547.3928
$ws.Range("M63").Value = -3138.5
$ws.Range("N63").Value = -4919.3928
$ws.Range("H66").Value = 3565.8667
$ws.Range("I66").Value = 3824.5
$ws.Range("J66").Value = 3547.3928
$ws.Range("K66").Value = 19122.5
$ws.Range("L66").Value = 17736.964
$ws.Range("M66").Value = -15690.5
$ws.Range("N66").Value = -24600.964

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 19062
$ws.Range("J82").Value = 53332
$ws.Range("L82").Value = 53332
$ws.Range("N82").Value = -54098
$ws.Range("H85").Value = 19062
$ws.Range("J85").Value = 53332
$ws.Range("L85").Value = 53332
$ws.Range("N85").Value = -55984

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 43.07143
$ws.Range("I7").Value = 53.842106
$ws.Range("J7").Value = 20.333334
$ws.Range("K7").Value = 53.842106
$ws.Range("L7").Value = 20.333334
$ws.Range("M7").Value = 59.157894
$ws.Range("N7").Value = -246.333334
$ws.Range("H16").Value = 1525.95
$ws.Range("J16").Value = 2312.8333
$ws.Range("L16").Value = 2312.8333
$ws.Range("N16").Value = -2886.8333
$ws.Range("H31").Value = 2142.125
$ws.Range("J31").Value = 3499.6667
$ws.Range("L31").Value = 3499.6667
$ws.Range("N31").Value = -4089.6667
$ws.Range("H34").Value = 2142.125
$ws.Range("J34").Value = 3499.6667
$ws.Range("L34").Value = 3499.6667
$ws.Range("N34").Value = -3903.6667
$ws.Range("H113").Value = 1525.95
$ws.Range("J113").Value = 2312.8333
$ws.Range("L113").Value = 2312.8333
$ws.Range("N113").Value = -6652.8333
$ws.Range("H134").Value = 4142.5
$ws.Range("I134").Value = 3896.389
$ws.Range("K134").Value = 11689.167
$ws.Range("M134").Value = -9154.167000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 890
$ws.Range("I8").Value = 890
$ws.Range("K8").Value = 2670
$ws.Range("M8").Value = -2531

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3500.8572
$ws.Range("I7").Value = 3252
$ws.Range("K7").Value = 3252
$ws.Range("M7").Value = -3140
$ws.Range("H126").Value = 3500.8572
$ws.Range("I126").Value = 3252
$ws.Range("K126").Value = 9756
$ws.Range("M126").Value = -7286

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 49992
$ws.Range("I40").Value = 49992
$ws.Range("K40").Value = 49992
$ws.Range("M40").Value = -49843
$ws.Range("H42").Value = 39000
$ws.Range("I42").Value = 39000
$ws.Range("K42").Value = 39000
$ws.Range("M42").Value = -38622
$ws.Range("H49").Value = 46496.332
$ws.Range("J49").Value = 49999
$ws.Range("L49").Value = 49999
$ws.Range("N49").Value = -50459
$ws.Range("H64").Value = 54555.5
$ws.Range("I64").Value = 49112
$ws.Range("J64").Value = 59999
$ws.Range("K64").Value = 49112
$ws.Range("L64").Value = 59999
$ws.Range("M64").Value = -48864
$ws.Range("N64").Value = -60495
$ws.Range("H67").Value = 54555.5
$ws.Range("I67").Value = 49112
$ws.Range("J67").Value = 59999
$ws.Range("K67").Value = 49112
$ws.Range("L67").Value = 59999
$ws.Range("M67").Value = -48254
$ws.Range("N67").Value = -61715
$ws.Range("H74").Value = 38484.4
$ws.Range("J74").Value = 36605.5
$ws.Range("L74").Value = 36605.5
$ws.Range("N74").Value = -38477.5
$ws.Range("H77").Value = 38484.4
$ws.Range("J77").Value = 36605.5
$ws.Range("L77").Value = 109816.5
$ws.Range("N77").Value = -119176.5
$ws.Range("H122").Value = 4473.193
$ws.Range("I122").Value = 2141.9778
$ws.Range("K122").Value = 6425.9334
$ws.Range("M122").Value = -3975.9334
